# Adds Family and Sample-related supporting sheets/classes to the ACR
# harmonized data model workbook, in line with current best practices:
#   - Adjusts "Sample" columns (splits FK/quantity fields, adds
#     aliquots + has_access_policy).
#   - Adds new "BiospecimenCollection" and "Aliquot" sheets.
#   - Adds new "Family", "FamilyRelationship" and "FamilyMember" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the existing "Sample" sheet header row / columns.
# ---------------------------------------------------------------------
$sample = $wb.Worksheets.Item("Sample")

$sample.Range("A1").Value = "biospecimen_collection"
$sample.Range("B1").Value = "parent_sample"
$sample.Range("C1").Value = "sample_type"
$sample.Range("D1").Value = "processing"
$sample.Range("E1").Value = "availablity_status"
$sample.Range("F1").Value = "storage_method"
$sample.Range("G1").Value = "quantity_number"
$sample.Range("H1").Value = "quantity_units"
$sample.Range("I1").Value = "aliquots"
$sample.Range("J1").Value = "has_access_policy"
$sample.Range("K1").Value = "id"
$sample.Range("L1").Value = "external_id"

# ---------------------------------------------------------------------
# 2. Insert "BiospecimenCollection" and "Aliquot" right after "Sample"
#    (and before "Subject").
# ---------------------------------------------------------------------
$biospecimenCollection = $wb.Worksheets.Add($null, $sample)
$biospecimenCollection.Name = "BiospecimenCollection"

$biospecimenCollection.Range("A1").Value = "age_at_collection"
$biospecimenCollection.Range("B1").Value = "method"
$biospecimenCollection.Range("C1").Value = "site"
$biospecimenCollection.Range("D1").Value = "spatial_qualifier"
$biospecimenCollection.Range("E1").Value = "laterality"
$biospecimenCollection.Range("F1").Value = "has_access_policy"
$biospecimenCollection.Range("G1").Value = "id"
$biospecimenCollection.Range("H1").Value = "external_id"

foreach ($col in @("B", "C", "D", "E")) {
    $validation = $biospecimenCollection.Range("$col`2:$col`1048576").Validation
    $validation.Add(3, 1, 1, '""')
    $validation.ShowInput = $false
    $validation.ShowError = $false
    $validation.InCellDropdown = $true
}

$aliquot = $wb.Worksheets.Add($null, $biospecimenCollection)
$aliquot.Name = "Aliquot"

$aliquot.Range("A1").Value = "availablity_status"
$aliquot.Range("B1").Value = "quantity_number"
$aliquot.Range("C1").Value = "quantity_units"
$aliquot.Range("D1").Value = "concentration_number"
$aliquot.Range("E1").Value = "concentration_unit"
$aliquot.Range("F1").Value = "has_access_policy"
$aliquot.Range("G1").Value = "id"
$aliquot.Range("H1").Value = "external_id"

$aliquotValidation = $aliquot.Range("A2:A1048576").Validation
$aliquotValidation.Add(3, 1, 1, '"available,unavailable"')
$aliquotValidation.ShowInput = $false
$aliquotValidation.ShowError = $false
$aliquotValidation.InCellDropdown = $true

# ---------------------------------------------------------------------
# 3. Append "Family", "FamilyRelationship" and "FamilyMember" at the
#    end of the workbook (after "AccessPolicy").
# ---------------------------------------------------------------------
$accessPolicy = $wb.Worksheets.Item("AccessPolicy")

$family = $wb.Worksheets.Add($null, $accessPolicy)
$family.Name = "Family"

$family.Range("A1").Value = "family_type"
$family.Range("B1").Value = "family_description"
$family.Range("C1").Value = "consanguinity"
$family.Range("D1").Value = "family_study_focus"
$family.Range("E1").Value = "family_members"
$family.Range("F1").Value = "family_relationships"
$family.Range("G1").Value = "has_access_policy"
$family.Range("H1").Value = "id"
$family.Range("I1").Value = "external_id"

$familyTypeValidation = $family.Range("A2:A1048576").Validation
$familyTypeValidation.Add(3, 1, 1, '"control_only,duo,proband_only,trio,trio_plus,other"')
$familyTypeValidation.ShowInput = $false
$familyTypeValidation.ShowError = $false
$familyTypeValidation.InCellDropdown = $true

$consanguinityValidation = $family.Range("C2:C1048576").Validation
$consanguinityValidation.Add(3, 1, 1, '"not_suspected,suspected,known_present,unknown"')
$consanguinityValidation.ShowInput = $false
$consanguinityValidation.ShowError = $false
$consanguinityValidation.InCellDropdown = $true

$familyRelationship = $wb.Worksheets.Add($null, $family)
$familyRelationship.Name = "FamilyRelationship"

$familyRelationship.Range("A1").Value = "family_member"
$familyRelationship.Range("B1").Value = "other_family_member"
$familyRelationship.Range("C1").Value = "relationship_code"
$familyRelationship.Range("D1").Value = "has_access_policy"
$familyRelationship.Range("E1").Value = "id"
$familyRelationship.Range("F1").Value = "external_id"

$familyMember = $wb.Worksheets.Add($null, $familyRelationship)
$familyMember.Name = "FamilyMember"

$familyMember.Range("A1").Value = "family_member"
$familyMember.Range("B1").Value = "family_role"
$familyMember.Range("C1").Value = "has_access_policy"
$familyMember.Range("D1").Value = "id"
$familyMember.Range("E1").Value = "external_id"
